$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I4").Value = 1.498352733344487
$ws.Range("J4").Value = 0.7413492207625757
$ws.Range("K4").Value = -0.7369451616989244
$ws.Range("L4").Value = 2.647587210590345
